# Update test data to current year (dates shifted forward by 366 days / 1 year).
$wb = $excel.ActiveWorkbook

# --- Bank In ---
$ws = $wb.Worksheets.Item("Bank In")
$ws.Range("A2").Value = 43952
$ws.Range("L2").Value = 43952
$ws.Range("A3").Value = 43953
$ws.Range("L3").Value = 43953
$ws.Range("A4").Value = 43954
$ws.Range("L4").Value = 43954
$ws.Range("A6").Value = 43982
$ws.Range("L2:L4").Select()

# --- Bank Out ---
$ws = $wb.Worksheets.Item("Bank Out")
$ws.Range("A2").Value = 43952
$ws.Range("L2").Value = 43952
$ws.Range("A3").Value = 43953
$ws.Range("L3").Value = 43953
$ws.Range("A4").Value = 43954
$ws.Range("L4").Value = 43954
$ws.Range("A5").Value = 43955
$ws.Range("L5").Value = 43955
$ws.Range("A6").Value = 43963
$ws.Range("L6").Value = 43963
$ws.Range("A7").Value = 43964
$ws.Range("L7").Value = 43964
$ws.Range("A9").Value = 43982
$ws.Range("A2:A9").Select()

# --- CredCard1 ---
$ws = $wb.Worksheets.Item("CredCard1")
$ws.Range("A2").Value = 43952
$ws.Range("H2").Value = 43952
$ws.Range("A3").Value = 43953
$ws.Range("H3").Value = 43953
$ws.Range("A4").Value = 43954
$ws.Range("H4").Value = 43954
$ws.Range("A6").Value = 43982
$ws.Range("H2:H4").Select()

# --- CredCard2 ---
$ws = $wb.Worksheets.Item("CredCard2")
$ws.Range("A2").Value = 43952
$ws.Range("H2").Value = 43952
$ws.Range("A3").Value = 43953
$ws.Range("H3").Value = 43953
$ws.Range("A4").Value = 43954
$ws.Range("H4").Value = 43954
$ws.Range("A6").Value = 43982
$ws.Range("H2:H4").Select()

# --- Expected In ---
$ws = $wb.Worksheets.Item("Expected In")
$ws.Range("A2").Value = 43952
$ws.Range("E2").Value = 43952
$ws.Range("A3").Value = 43953
$ws.Range("E3").Value = 43953
$ws.Range("A4").Value = 43954
$ws.Range("E4").Value = 43954
$ws.Range("A6").Value = 43983
$ws.Range("A7").Value = 43984
$ws.Range("A8").Value = 43985
$ws.Range("E2:E4").Select()

# --- CredCard3 ---
$ws = $wb.Worksheets.Item("CredCard3")
$ws.Range("A2").Value = 43952
$ws.Range("A3").Value = 43953
$ws.Range("A4").Value = 43954
$ws.Range("A5").Value = 43983
$ws.Range("A2:A5").Select()

# --- Savings (values only; selection unchanged) ---
$ws = $wb.Worksheets.Item("Savings")
$ws.Range("A2").Value = 43891
$ws.Range("A3").Value = 43922
$ws.Range("A4").Value = 43952

# --- Budget Out becomes the active sheet/tab (selection on it is unchanged: G18) ---
$ws = $wb.Worksheets.Item("Budget Out")
$ws.Activate()

Write-Host "edit complete"
